# QUBES Code First commit
# Re-create the edits made to TestData.xlsx:
#  - A4 is reset to a plain (unstyled) cell containing "NewData 98"
#  - C5 is updated to "Var1-VS1P320220008"
#  - C6 is updated to "SBH035177"
#  - The sheet's scrolled view is nudged one column to the right (F1 -> G1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 previously carried a border style inherited from the row above; the
# updated value lives in a bare cell, so clear formatting before writing it.
$ws.Range("A4").Clear()
$ws.Range("A4").Value = "NewData 98"

$ws.Range("C5").Value = "Var1-VS1P320220008"
$ws.Range("C6").Value = "SBH035177"

# Scroll the visible window one column to the right so the top-left visible
# cell moves from F1 to G1. (The current selection, L2, is left untouched.)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 7
